# Update the "РЕФЕРАТ" (abstract) page of the КП document:
#   1. Replace the (misspelled) project-title phrase inside the
#      bibliographic-description paragraph with the corrected title.
#   2. Bump the total page count 47 -> 52 (two places: the
#      bibliographic description and the "Пояснительная записка" line).
#   3. Bump the appendix count 2 -> 3.
#   4. Replace the running/caps title paragraph with the corrected,
#      fully upper-cased title.
#
# wdReplaceOne = 1 (replace just the first/next match) is used
# throughout so that runs of identical text elsewhere in the document
# are left untouched; wdFindContinue = 1 for Wrap.

$d = $word.ActiveDocument

# 1) Bibliographic-description paragraph: fix the project title, keep
#    the rest (" / Д.Л. Богомаз. - Минск : БГУИР, 20") intact.
$d.Content.Find.Execute(
    " Програмное средство хранения и одновремнного редактирования заметок / Д.Л. Богомаз. – Минск : БГУИР, 20",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Программное средство для реализации менеджера заметок и задач с возможностью визуализации и структуризации материала / Д.Л. Богомаз. – Минск : БГУИР, 20",
    1) | Out-Null

# 2) Page count 47 -> 52 in the bibliographic description.
$d.Content.Find.Execute(
    ". – 47 с.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ". – 52 с.",
    1) | Out-Null

# 3) Page count 47 -> 52 in "Пояснительная записка ...".
$d.Content.Find.Execute(
    "Пояснительная записка 47 с.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Пояснительная записка 52 с.",
    1) | Out-Null

# 4) Appendix count 2 -> 3.
$d.Content.Find.Execute(
    "источника, 2 приложения",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "источника, 3 приложения",
    1) | Out-Null

# 5) Running caps-styled title paragraph. MatchCase is turned off here
#    because this run carries the w:caps formatting, so search case
#    doesn't need to (and must not) be forced to match the lower-case
#    text stored in the run.
$d.Content.Find.Execute(
    "Програмное средство хранения и одновремнного редактирования заметок ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "ПРОГРАММНОЕ СРЕДСТВО ДЛЯ РЕАЛИЗАЦИИ МЕНЕДЖЕРА ЗАМЕТОК И ЗАДАЧ С ВОЗМОЖНОСТЬЮ ВИЗУАЛИЗАЦИИ И СТРУКТУРИЗАЦИИ МАТЕРИАЛА",
    1) | Out-Null
